# Update NATmi TPM-derived edge-expression metrics for the Fgf15-Fgfr1
# ligand/receptor table (Sheet1). Only the numeric columns F:T (ligand/
# receptor detection, expression and specificity figures) were recomputed
# with the new TPM values; the label columns (A sending cluster, B ligand,
# C receptor, D target cluster) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("G2").Value  = 0.018508
$ws.Range("H2").Value  = 0.055524
$ws.Range("I2").Value  = 0.2347866901774728
$ws.Range("J2").Value  = 0.3151800006811757
$ws.Range("M2").Value  = 5.978421000000001
$ws.Range("N2").Value  = 17.935263
$ws.Range("O2").Value  = 0.05704457007880161
$ws.Range("P2").Value  = 0.06242884486533885
$ws.Range("Q2").Value  = 0.110648615868
$ws.Range("R2").Value  = 0.9958375428120001
$ws.Range("S2").Value  = 0.01339330580139873
$ws.Range("T2").Value  = 0.01967632336718251

# Row 3 (Target cluster: FAPs)
$ws.Range("G3").Value  = 0.018508
$ws.Range("H3").Value  = 0.055524
$ws.Range("I3").Value  = 0.2347866901774728
$ws.Range("J3").Value  = 0.3151800006811757
$ws.Range("O3").Value  = 0.6646576013185088
$ws.Range("P3").Value  = 0.7273927426214574
$ws.Range("Q3").Value  = 1.28922776542
$ws.Range("R3").Value  = 11.60304988878
$ws.Range("S3").Value  = 0.1560527583148709
$ws.Range("T3").Value  = 0.2292596451149132

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("G4").Value  = 0.018508
$ws.Range("H4").Value  = 0.055524
$ws.Range("I4").Value  = 0.2347866901774728
$ws.Range("J4").Value  = 0.3151800006811757
$ws.Range("M4").Value  = 1.290243
$ws.Range("N4").Value  = 3.870729
$ws.Range("O4").Value  = 0.01231116999491725
$ws.Range("P4").Value  = 0.01347318632889677
$ws.Range("Q4").Value  = 0.023879817444
$ws.Range("R4").Value  = 0.214918356996
$ws.Range("S4").Value  = 0.002890498855318837
$ws.Range("T4").Value  = 0.004246478876319292

# Row 5 (Target cluster: MuSCs)
$ws.Range("G5").Value  = 0.018508
$ws.Range("H5").Value  = 0.055524
$ws.Range("I5").Value  = 0.2347866901774728
$ws.Range("J5").Value  = 0.3151800006811757
$ws.Range("M5").Value  = 27.1166075
$ws.Range("N5").Value  = 54.233215
$ws.Range("O5").Value  = 0.2587397603536297
$ws.Range("P5").Value  = 0.1887743138075849
$ws.Range("Q5").Value  = 0.50187417161
$ws.Range("R5").Value  = 3.01124502966
$ws.Range("S5").Value  = 0.0607486519507412
$ws.Range("T5").Value  = 0.0594978883544631

# Row 6 (Target cluster: Resolving-Mac)
$ws.Range("G6").Value  = 0.018508
$ws.Range("H6").Value  = 0.055524
$ws.Range("I6").Value  = 0.2347866901774728
$ws.Range("J6").Value  = 0.3151800006811757
$ws.Range("M6").Value  = 0.759494
$ws.Range("N6").Value  = 2.278482
$ws.Range("O6").Value  = 0.00724689825414258
$ws.Range("P6").Value  = 0.007930912376722157
$ws.Range("Q6").Value  = 0.014056714952
$ws.Range("R6").Value  = 0.126510434568
$ws.Range("S6").Value  = 0.001701475255143042
$ws.Range("T6").Value  = 0.002499664968297634

# Row 7 (Sending cluster: Resolving-Mac; Target cluster: ECs)
$ws.Range("F7").Value  = 0.5
$ws.Range("G7").Value  = 0.060321
$ws.Range("H7").Value  = 0.120642
$ws.Range("I7").Value  = 0.7652133098225272
$ws.Range("J7").Value  = 0.6848199993188243
$ws.Range("M7").Value  = 5.978421000000001
$ws.Range("N7").Value  = 17.935263
$ws.Range("O7").Value  = 0.05704457007880161
$ws.Range("P7").Value  = 0.06242884486533885
$ws.Range("Q7").Value  = 0.3606243331410001
$ws.Range("R7").Value  = 2.163745998846
$ws.Range("S7").Value  = 0.04365126427740288
$ws.Range("T7").Value  = 0.04275252149815634

# Row 8 (Sending cluster: Resolving-Mac; Target cluster: FAPs)
$ws.Range("F8").Value  = 0.5
$ws.Range("G8").Value  = 0.060321
$ws.Range("H8").Value  = 0.120642
$ws.Range("I8").Value  = 0.7652133098225272
$ws.Range("J8").Value  = 0.6848199993188243
$ws.Range("O8").Value  = 0.6646576013185088
$ws.Range("P8").Value  = 0.7273927426214574
$ws.Range("Q8").Value  = 4.201832074665
$ws.Range("R8").Value  = 25.21099244799
$ws.Range("S8").Value  = 0.5086048430036378
$ws.Range("T8").Value  = 0.4981330975065442

# Row 9 (Sending cluster: Resolving-Mac; Target cluster: Inflammatory-Mac)
$ws.Range("F9").Value  = 0.5
$ws.Range("G9").Value  = 0.060321
$ws.Range("H9").Value  = 0.120642
$ws.Range("I9").Value  = 0.7652133098225272
$ws.Range("J9").Value  = 0.6848199993188243
$ws.Range("M9").Value  = 1.290243
$ws.Range("N9").Value  = 3.870729
$ws.Range("O9").Value  = 0.01231116999491725
$ws.Range("P9").Value  = 0.01347318632889677
$ws.Range("Q9").Value  = 0.07782874800300001
$ws.Range("R9").Value  = 0.466972488018
$ws.Range("S9").Value  = 0.009420671139598418
$ws.Range("T9").Value  = 0.009226707452577481

# Row 10 (Sending cluster: Resolving-Mac; Target cluster: MuSCs)
$ws.Range("F10").Value = 0.5
$ws.Range("G10").Value = 0.060321
$ws.Range("H10").Value = 0.120642
$ws.Range("I10").Value = 0.7652133098225272
$ws.Range("J10").Value = 0.6848199993188243
$ws.Range("M10").Value = 27.1166075
$ws.Range("N10").Value = 54.233215
$ws.Range("O10").Value = 0.2587397603536297
$ws.Range("P10").Value = 0.1887743138075849
$ws.Range("Q10").Value = 1.6357008810075
$ws.Range("R10").Value = 6.54280352403
$ws.Range("S10").Value = 0.1979911084028885
$ws.Range("T10").Value = 0.1292764254531218

# Row 11 (Sending cluster: Resolving-Mac; Target cluster: Resolving-Mac)
$ws.Range("F11").Value = 0.5
$ws.Range("G11").Value = 0.060321
$ws.Range("H11").Value = 0.120642
$ws.Range("I11").Value = 0.7652133098225272
$ws.Range("J11").Value = 0.6848199993188243
$ws.Range("M11").Value = 0.759494
$ws.Range("N11").Value = 2.278482
$ws.Range("O11").Value = 0.00724689825414258
$ws.Range("P11").Value = 0.007930912376722157
$ws.Range("Q11").Value = 0.045813437574
$ws.Range("R11").Value = 0.274880625444
$ws.Range("S11").Value = 0.005545422998999538
$ws.Range("T11").Value = 0.005431247408424523
